# Re-running the backward-elimination export re-stamped every OLS summary
# block with the new run's Date/Time line (the underlying numbers are
# unchanged - this mirrors wrapping the Excel write in a try/except and
# re-running it, which is all the commit actually changed in the data).

$wb = $excel.ActiveWorkbook

# Sheets are in the same order as the OOXML parts (sheet1.xml .. sheet29.xml):
# tab "46" down to tab "18". Every sheet has its statsmodels OLS summary
# text in cell B2. The first 21 reruns landed at 23:19:04, the last 8 at
# 23:19:05 (matching the one-second tick crossed while the export looped).
$newTimes = @(
    "23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04",
    "23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04",
    "23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04","23:19:04",
    "23:19:05","23:19:05","23:19:05","23:19:05","23:19:05","23:19:05","23:19:05",
    "23:19:05"
)

$oldDate = "Sun, 29 Dec 2019"
$newDate = "Wed, 01 Jan 2020"
$oldTime = "16:11:27"
$oldTime2 = "16:11:28"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Cells.Item(2, 2)
    $text = $cell.Value2

    if ($text -ne $null -and $text.Contains("Date:") -and $text.Contains("Time:")) {
        $newTime = $newTimes[$i - 1]
        $updated = $text.Replace($oldDate, $newDate)
        $updated = $updated.Replace($oldTime, $newTime)
        $updated = $updated.Replace($oldTime2, $newTime)
        $cell.Value = $updated
    }
}
